$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Google'
$ws.Range('B1').Value = 'Amazon'
$ws.Range('A2').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/robertprescott&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjACegQICBAB&usg=AOvVaw2FkHodd__nRdVJRuAO1jr5'
$ws.Range('B2').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/kyleroysmith&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAAegQIBRAB&usg=AOvVaw2emxTBlugLKjqgZAd6mefv'
$ws.Range('A3').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/fernandoterrero&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjADegQIBxAB&usg=AOvVaw3AhEzheb6pqW58ETymWM_t'
$ws.Range('B3').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/amymil&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjABegQICRAB&usg=AOvVaw1MYjGwuziqtwzdpWlgz7v6'
$ws.Range('A4').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/petersenjordan&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjAEegQICRAB&usg=AOvVaw2zkQOovCLCquzUxvZngMz2'
$ws.Range('B4').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/tripti-vishwakarma&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjACegQIBxAB&usg=AOvVaw1Kv2GnYCmAR_vCbYRfQb9h'
$ws.Range('A5').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/jacq-morrill-12a41983&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjAFegQIABAB&usg=AOvVaw0pFErZl4HLXuQiNz3yYwxW'
$ws.Range('B5').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/emucken&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjADegQICBAB&usg=AOvVaw1BKb4_hUDIufiwCZwKXlFn'
$ws.Range('A6').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/ericbarradas&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjAHegQIAxAB&usg=AOvVaw0pBdsFtdBrCX5EcgkEg2sf'
$ws.Range('B6').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/jennifer-kemp-71057251&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAEegQIChAB&usg=AOvVaw1dwP37JjfUA5rKm8MJI8Eq'
$ws.Range('A7').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/latasha-smith-73996b4b&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjAIegQIARAB&usg=AOvVaw0yhKNcntbMO2K1Tb3URlNz'
$ws.Range('B7').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/saramorgan89&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAFegQIAxAB&usg=AOvVaw0tiQke8z8I1Yc34lzfBPqc'
$ws.Range('A8').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/chantal-penton-1165b461&sa=U&ved=2ahUKEwinppSc0cHwAhWDup4KHVatAA8QFjAJegQIBRAB&usg=AOvVaw08Dgm_S3q7ZVhBOAx2ixfH'
$ws.Range('B8').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/danadworman&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAGegQIAhAB&usg=AOvVaw3JFQGyxV-_XFMaFeVRl1fG'
$ws.Range('B9').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/claireoc&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAHegQIBBAB&usg=AOvVaw1aqLKDBcmS6Jp2qkDSOmMr'
$ws.Range('B10').Value = 'https://www.google.com/url?q=https://www.linkedin.com/in/joshmayaws&sa=U&ved=2ahUKEwiHk4at0cHwAhVPrJ4KHRc2C3wQFjAJegQIABAB&usg=AOvVaw3AWcFhLaTzgYMXkj1nms59'

# Set column widths for A and B to 200 (matches ColumnWidth 199.16666666666666 -> OOXML width 200)
$ws.Columns.Item(1).ColumnWidth = 199.16666666666666
$ws.Columns.Item(2).ColumnWidth = 199.16666666666666
